$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells stay text, matching original inlineStr formatting,
# so Excel does not auto-convert numeric-looking strings (e.g. "1.003") into numbers.
$dCells = @("D2","D3","D4","D5","D6","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated coin data values
$ws.Range("D2").Value = "25.822.97"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "1.635.09"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "215.08"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").Value = "0.5083"
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("D8").Value = "0.2578"
$ws.Range("E8").Value = "  +0.57%  "
$ws.Range("D9").Value = "0.06427"
$ws.Range("E9").Value = "  +1.49%  "
$ws.Range("D10").Value = "20.25"
$ws.Range("D11").Value = "0.07792"
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").Value = "1.649.02"
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("D13").Value = "4.248"
$ws.Range("E13").Value = "  +0.06%  "
$ws.Range("D14").Value = "1.859.98"
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("D15").Value = "0.5585"
$ws.Range("E15").Value = "  +1.14%  "
$ws.Range("D16").Value = "0.0₅7651"
$ws.Range("E16").Value = "  +0.85%  "
$ws.Range("D17").Value = "63.20"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("D18").Value = "25.831.90"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").Value = "4.362"
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("D21").Value = "191.88"
$ws.Range("E21").Value = "  -1.32%  "
$ws.Range("D22").Value = "9.916"
$ws.Range("E22").Value = "  +0.64%  "
$ws.Range("D23").Value = "6.140"
$ws.Range("E23").Value = "  +2.08%  "
$ws.Range("D24").Value = "1.004"
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("D25").Value = "1.762"
$ws.Range("E25").Value = "  -6.66%  "
$ws.Range("D26").Value = "139.01"
$ws.Range("E26").Value = "  -2.21%  "
$ws.Range("D27").Value = "0.1228"
$ws.Range("E27").Value = "  -2.24%  "
$ws.Range("D28").Value = "6.817"
$ws.Range("E28").Value = "  +0.86%  "
$ws.Range("D29").Value = "15.56"
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("D30").Value = "1.241"
$ws.Range("D31").Value = "0.04946"
$ws.Range("E31").Value = "  +0.94%  "
$ws.Range("D32").Value = "3.295"
$ws.Range("E32").Value = "  +1.87%  "
$ws.Range("D33").Value = "3.250"
$ws.Range("E33").Value = "  +2.25%  "
$ws.Range("D34").Value = "1.567"
$ws.Range("E34").Value = "  +1.29%  "
$ws.Range("D35").Value = "2.387"
$ws.Range("E35").Value = "  +0.65%  "
$ws.Range("D36").Value = "0.9001"
$ws.Range("E36").Value = "  +0.50%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "0.5572"
$ws.Range("E37").Value = "  +0.92%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "2.566"
$ws.Range("E38").Value = "  +1.16%  "
$ws.Range("D39").Value = "1.130.08"
$ws.Range("E39").Value = "  +1.37%  "
$ws.Range("D40").Value = "0.01569"
$ws.Range("E40").Value = "  +1.02%  "
$ws.Range("D41").Value = "0.9964"
$ws.Range("E41").Value = "  -0.51%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "98.99"
$ws.Range("E42").Value = "  +1.39%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "5.449"
$ws.Range("E43").Value = "  -2.20%  "
$ws.Range("D44").Value = "0.7986"
$ws.Range("E44").Value = "  +0.48%  "
$ws.Range("D45").Value = "0.0₈113"
$ws.Range("E45").Value = "  -3.16%  "
$ws.Range("D46").Value = "55.55"
$ws.Range("E46").Value = "  +1.53%  "
$ws.Range("D47").Value = "0.4260"
$ws.Range("E47").Value = "  -3.92%  "
$ws.Range("D48").Value = "7.782"
$ws.Range("E48").Value = "  +2.84%  "
$ws.Range("D49").Value = "0.05026"
$ws.Range("E49").Value = "  -2.11%  "
$ws.Range("D50").Value = "0.9969"
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("E51").Value = "  +0.47%  "
